$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its original text formatting (Excel would
# otherwise auto-coerce values like "1.002" or "0.07660" into numbers).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.468.97"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.852.58"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.46"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6299"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07660"
$ws.Range("E8").Value = "  +0.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.86"
$ws.Range("E10").Value = "  +1.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "2.140.37"
$ws.Range("E11").Value = "  +15.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07752"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.041"
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6821"
$ws.Range("E14").Value = "  +0.51%  "
$ws.Range("E15").Value = "  -5.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.52"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.200"
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.596.57"
$ws.Range("E18").Value = "  +0.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "229.23"
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.36"
$ws.Range("E20").Value = "  -0.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.474"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "157.70"
$ws.Range("E24").Value = "  +0.28%  "
$ws.Range("E25").Value = "  -0.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.439"
$ws.Range("E26").Value = "  +1.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.75"
$ws.Range("E27").Value = "  +0.76%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.394"
$ws.Range("E28").Value = "  +7.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.465"
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05610"
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.138"
$ws.Range("E31").Value = "  +0.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.068"
$ws.Range("E32").Value = "  +0.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.848"
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.167"
$ws.Range("E34").Value = "  +0.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7012"
$ws.Range("E35").Value = "  -1.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.592"
$ws.Range("E36").Value = "  +0.41%  "
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.231.72"
$ws.Range("E38").Value = "  -0.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.744"
$ws.Range("E39").Value = "  -1.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.468"
$ws.Range("E40").Value = "  +0.86%  "
$ws.Range("E41").Value = "  +0.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.002"
$ws.Range("E42").Value = "  +0.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "102.74"
$ws.Range("E43").Value = "  +0.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "66.18"
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.208"
$ws.Range("E45").Value = "  +0.84%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000117"
$ws.Range("E46").Value = "  -2.62%  "
$ws.Range("B47").Value = "TheSandbox"
$ws.Range("C47").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4032"
$ws.Range("E47").Value = "  +0.42%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.033"
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1155"
$ws.Range("E49").Value = "  +3.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.683"
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05715"
$ws.Range("E51").Value = "  +0.10%  "
